# Update for revisions post peer-review.
# Target sheet is "Trends".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trends")
$ws.Activate()

# --- Row 2 / Row 3: small value tweaks, formatting untouched ---
$ws.Range("C2").Value = 300
$ws.Range("E2").Value = 310
$ws.Range("C3").Value = 300
$ws.Range("E3").Value = 310

# --- Row 4: whole row of data replaced; cells lose their explicit
#     (bordered) formatting and become plain/default-styled, and the
#     AE4 shared formula is replaced by a literal value. ---
$ws.Range("B4:AE4").Style = "Normal"

$ws.Range("B4").Value = 2600
$ws.Range("C4").Value = 285
$ws.Range("D4").Value = 2375
$ws.Range("E4").Value = 295
$ws.Range("F4").Value = 5000
$ws.Range("G4").Value = 11750
$ws.Range("H4").Value = 16
$ws.Range("I4").Value = 35
$ws.Range("J4").Value = 7.5
$ws.Range("K4").Value = 0.5
$ws.Range("L4").Value = 17.5
$ws.Range("M4").Value = 5
$ws.Range("N4").Value = -1.5
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 0.625
$ws.Range("Q4").Value = 1.2
$ws.Range("R4").Value = 140
$ws.Range("S4").Value = 27.5
$ws.Range("T4").Value = 41
$ws.Range("U4").Value = 10.5
$ws.Range("V4").Value = 0.525
$ws.Range("W4").Value = 0.1
$ws.Range("X4").Value = 18
$ws.Range("Y4").Value = 42
$ws.Range("Z4").Value = 3.5
$ws.Range("AA4").Value = 1.25
$ws.Range("AB4").Value = 1000
$ws.Range("AC4").Value = 150000
$ws.Range("AD4").Value = 600
$ws.Range("AE4").Value = 0.36

# --- Row 5: whole row of data replaced too. Column A keeps its
#     original style (s=10), everything else becomes plain/default,
#     the AE5 shared formula becomes a literal value, and the row
#     itself drops its explicit height / thick-bottom-border flag. ---
$ws.Range("B5:AE5").Style = "Normal"

$ws.Range("B5").Value = 2625
$ws.Range("C5").Value = 325
$ws.Range("D5").Value = 2325
$ws.Range("E5").Value = 320
$ws.Range("F5").Value = 6000
$ws.Range("G5").Value = 14500
$ws.Range("H5").Value = 22.5
$ws.Range("I5").Value = 30
$ws.Range("J5").Value = 7
$ws.Range("K5").Value = 1.85
$ws.Range("L5").Value = 21
$ws.Range("M5").Value = 6.5
$ws.Range("N5").Value = -1.75
$ws.Range("O5").Value = 1.75
$ws.Range("P5").Value = 0.8300000000000001
$ws.Range("Q5").Value = 1.05
$ws.Range("R5").Value = 120
$ws.Range("S5").Value = 32.5
$ws.Range("T5").Value = 43.5
$ws.Range("U5").Value = 11.5
$ws.Range("V5").Value = 0.675
$ws.Range("W5").Value = 0.14
$ws.Range("X5").Value = 18
$ws.Range("Y5").Value = 55
$ws.Range("Z5").Value = 3.75
$ws.Range("AA5").Value = 1.75
$ws.Range("AB5").Value = 1000
$ws.Range("AC5").Value = 150000
$ws.Range("AD5").Value = 600
$ws.Range("AE5").Value = 0.36

# Row 5 no longer has an explicit/custom row height or thick bottom
# border once its content was replaced - let Excel recompute it.
$ws.Rows.Item(5).AutoFit()

# --- View state: scrolled over one column further right, and the
#     active selection moved from W2 to Y5. ---
$excel.ActiveWindow.ScrollColumn = 16
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("Y5").Select()
